$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.788.95'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.798.68'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.82%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.82%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.43'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.67%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4574'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.53%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3723'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07336'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8564'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.38'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.814.96'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.600'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.343'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07074'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.06'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.007'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.87%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008636'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.75'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.829.25'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.283'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.79'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.038.65'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.906'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.25'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.185'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.42'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.226'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.68'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08871'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7649'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.161'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.471'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.898'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.59%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.122'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.44%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01954'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05209'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.42%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.189'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.79%  '

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.885'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5242'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.311'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +6.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1663'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.522'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4970'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.59%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.30'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.15%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.88'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.76%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.657'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.37%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06308'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.52%  '
